$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price column (D) keeps its original text representation instead of
# being auto-coerced to a Number by the .Value setter (mirrors the source data,
# which stores prices as literal text, e.g. "26.086.61", "1.666.40", "1.002").
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '26.086.61'
$ws.Range("E2").Value = '  -0.71%  '

$ws.Range("D3").Value = '1.666.40'
$ws.Range("E3").Value = '  -1.32%  '

$ws.Range("D4").Value = '1.002'
$ws.Range("E4").Value = '  -0.59%  '

$ws.Range("D5").Value = '209.57'
$ws.Range("E5").Value = '  -3.61%  '

$ws.Range("E6").Value = '  -1.60%  '

$ws.Range("E7").Value = '  -0.57%  '

$ws.Range("E8").Value = '  -3.23%  '

$ws.Range("D9").Value = '0.06287'
$ws.Range("E9").Value = '  -1.90%  '

$ws.Range("D11").Value = '0.07538'
$ws.Range("E11").Value = '  -1.68%  '

$ws.Range("D12").Value = '1.670.19'
$ws.Range("E12").Value = '  -1.23%  '

$ws.Range("D13").Value = '4.437'
$ws.Range("E13").Value = '  -1.92%  '

$ws.Range("D15").Value = '66.55'
$ws.Range("E15").Value = '  +0.10%  '

$ws.Range("D16").Value = '0.000007933'
$ws.Range("E16").Value = '  -4.89%  '

$ws.Range("D17").Value = '26.137.00'
$ws.Range("E17").Value = '  -0.66%  '

$ws.Range("E18").Value = '  -0.56%  '

$ws.Range("D19").Value = '4.718'
$ws.Range("E19").Value = '  -3.24%  '

$ws.Range("D20").Value = '186.24'
$ws.Range("E20").Value = '  -2.32%  '

$ws.Range("E21").Value = '  -4.95%  '

$ws.Range("D22").Value = '6.157'
$ws.Range("E22").Value = '  -1.34%  '

$ws.Range("E23").Value = '  -0.57%  '

$ws.Range("D24").Value = '149.77'
$ws.Range("E24").Value = '  +0.81%  '

$ws.Range("D25").Value = '0.1248'
$ws.Range("E25").Value = '  -2.99%  '

$ws.Range("D26").Value = '7.479'
$ws.Range("E26").Value = '  -4.63%  '

$ws.Range("D27").Value = '15.88'
$ws.Range("E27").Value = '  +0.68%  '

$ws.Range("D28").Value = '0.06308'
$ws.Range("E28").Value = '  +1.91%  '

$ws.Range("D29").Value = '1.356'
$ws.Range("E29").Value = '  -1.45%  '

$ws.Range("D30").Value = '1.273'
$ws.Range("E30").Value = '  -3.89%  '

$ws.Range("E31").Value = '  -2.70%  '

$ws.Range("D32").Value = '3.409'
$ws.Range("E32").Value = '  -4.73%  '

$ws.Range("D33").Value = '1.636'
$ws.Range("E33").Value = '  -2.44%  '

$ws.Range("D34").Value = '0.9958'
$ws.Range("E34").Value = '  -3.00%  '

$ws.Range("D35").Value = '0.6023'
$ws.Range("E35").Value = '  -3.14%  '

$ws.Range("E36").Value = '  -0.70%  '

$ws.Range("D37").Value = '2.728'
$ws.Range("E37").Value = '  -1.24%  '

$ws.Range("D38").Value = '1.109.14'
$ws.Range("E38").Value = '  +0.25%  '

$ws.Range("E39").Value = '  -0.39%  '

$ws.Range("D40").Value = '0.01613'
$ws.Range("E40").Value = '  -1.99%  '

$ws.Range("D41").Value = '0.8727'
$ws.Range("E41").Value = '  -0.89%  '

$ws.Range("E43").Value = '  -1.06%  '

$ws.Range("D44").Value = '1.819.25'
$ws.Range("E44").Value = '  -1.11%  '

$ws.Range("D45").Value = '0.00000000108'
$ws.Range("E45").Value = '  -3.24%  '

$ws.Range("D46").Value = '55.38'
$ws.Range("E46").Value = '  -3.80%  '

$ws.Range("D47").Value = '0.9996'
$ws.Range("E47").Value = '  -0.34%  '

$ws.Range("D48").Value = '8.012'
$ws.Range("E48").Value = '  -1.66%  '

$ws.Range("E49").Value = '  -0.83%  '

$ws.Range("D50").Value = '0.4243'
$ws.Range("E50").Value = '  -1.30%  '

$ws.Range("D51").Value = '5.958'
$ws.Range("E51").Value = '  -1.63%  '
